$wb = $excel.ActiveWorkbook

# ==========================================================================
# 1. "Logs" sheet: append two new mail-log rows (13 and 14)
# ==========================================================================
$wsLogs = $wb.Worksheets.Item("Logs")

# --- Row 13: Undelivered Mail Returned to Sender ---
$wsLogs.Cells.Item(13,1).Value = "Undelivered Mail Returned to Sender"
$wsLogs.Cells.Item(13,2).Value = "mailer-daemon@mail.zoho.eu"
$cell_C13 = @"
This message was created automatically by mail delivery software.
 A message that you sent could not be delivered to one or more of its recipients. This is a permanent error. 
support@testbedrijf123.nl, ERROR CODE :550 - 5.0.0 Invalid Recipients.
"@
$wsLogs.Cells.Item(13,3).Value = $cell_C13
$wsLogs.Cells.Item(13,4).Value = "Overig"
$cell_E13 = @"
Geachte afzender,
Hartelijk dank voor uw bericht. 
Het lijkt erop dat uw e-mail niet correct is afgeleverd bij onze support@testbedrijf123.nl mailbox vanwege een fout met betrekking tot ongeldige ontvangers (ERROR CODE: 550 - 5.0.0). 
Om dit probleem te verhelpen, verzoeken wij u vriendelijk om de e-mail opnieuw te verzenden en ervoor te zorgen dat het e-mailadres support@testbedrijf123.nl correct is ingevoerd in de ontvangerlijst. Mocht u verdere hulp nodig hebben, aarzel dan niet om contact met ons op te nemen.
Met vriendelijke groet,
[Naam]
E-mailassistent bij Testbedrijf123
"@
$wsLogs.Cells.Item(13,5).Value = $cell_E13
$wsLogs.Cells.Item(13,6).Value = "2025-08-14 20:43:57"
$wsLogs.Cells.Item(13,7).Value = "Ja"
$wsLogs.Cells.Item(13,8).Value = "Nee"
$wsLogs.Cells.Item(13,9).Value = "Ja"
$wsLogs.Cells.Item(13,10).Value = "Nee"

# Row 13 holds multi-line content (C13/E13); re-fit its height back to the
# sheet default instead of leaving Excel's auto-expanded custom row height.
$wsLogs.Rows.Item(13).AutoFit()

# --- Row 14: CE-certificaten verzoek ---
$wsLogs.Cells.Item(14,1).Value = "CE-certificaten verzoek"
$wsLogs.Cells.Item(14,2).Value = "inkoop@testbedrijf123.nl"
$wsLogs.Cells.Item(14,3).Value = "Kun je mij de CE-certificaten van de EcoPro-700 sturen?"
$wsLogs.Cells.Item(14,4).Value = "Kwaliteit / Certificaten"
$wsLogs.Cells.Item(14,5).Value = "Bedankt, we hebben dit doorgestuurd naar kwaliteit@testbedrijf123.nl."
$wsLogs.Cells.Item(14,6).Value = "2025-08-14 20:44:07"
$wsLogs.Cells.Item(14,7).Value = "Nee"
$wsLogs.Cells.Item(14,8).Value = "Ja"
$wsLogs.Cells.Item(14,9).Value = "Nee"
$wsLogs.Cells.Item(14,10).Value = "Nee"

# --------------------------------------------------------------------------
# Extend the conditional-formatting ranges (D/G/H/I/J 2:12 -> 2:14) so the
# two new rows get the same highlighting rules as the rest of the table.
# --------------------------------------------------------------------------
$wsLogs.Range("D2:D12").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("D2:D14"))
$wsLogs.Range("G2:G12").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("G2:G14"))
$wsLogs.Range("H2:H12").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("H2:H14"))
$wsLogs.Range("I2:I12").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("I2:I14"))
$wsLogs.Range("J2:J12").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("J2:J14"))

# ==========================================================================
# 2. "Dashboard" sheet: re-aggregate the category counts now that the two
#    new rows changed the tallies, adding a "Kwaliteit / Certificaten" row.
# ==========================================================================
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Cells.Item(3,1).Value = "Overig"
$wsDash.Cells.Item(3,2).Value = 2
$wsDash.Cells.Item(4,1).Value = "Planning / Afspraak"
$wsDash.Cells.Item(4,2).Value = 1
$wsDash.Cells.Item(5,1).Value = "Inkoop / Bestellingen"
$wsDash.Cells.Item(5,2).Value = 1
$wsDash.Cells.Item(6,1).Value = "Kwaliteit / Certificaten"
$wsDash.Cells.Item(6,2).Value = 1

# --------------------------------------------------------------------------
# 3. Chart on the Dashboard sheet: extend the category / value series ranges
#    from row 5 to row 6 to include the new "Kwaliteit / Certificaten" bar.
# --------------------------------------------------------------------------
$chart = $wsDash.ChartObjects(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!B1,'Dashboard'!`$A`$2:`$A`$6,'Dashboard'!`$B`$2:`$B`$6,1)"
